# Insert a new data row at row 77 (pushing existing rows 77..92 down to 78..93),
# and populate the new row 77 with its own values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(77).Insert()

$ws.Cells.Item(77, 1).Value = 4
$ws.Cells.Item(77, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(77, 3).Value = "Los Lagos"
$ws.Cells.Item(77, 4).Value = 45211
$ws.Cells.Item(77, 5).Value = 10
$ws.Cells.Item(77, 6).Value = 100112012
$ws.Cells.Item(77, 7).Value = "Espinaca"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 25
$ws.Cells.Item(77, 11).Value = 13000
$ws.Cells.Item(77, 12).Value = 13000
$ws.Cells.Item(77, 13).Value = 13000
$ws.Cells.Item(77, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(77, 15).Value = "Región Metropolitana"
$ws.Cells.Item(77, 16).Value = 1300
$ws.Cells.Item(77, 17).Value = 10
$ws.Cells.Item(77, 18).Value = "Hortaliza"
